# Refresh crypto price/volume figures (cols D and E) to match the latest
# GitHub Actions data pull, cell by cell, exactly as captured in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.059.04"
$ws.Range("E2").Value = "  +2.26%  "

$ws.Range("D3").Value = "2.301.76"
$ws.Range("E3").Value = "  +2.05%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'302.37"
$ws.Range("E5").Value = "  +1.43%  "

$ws.Range("D6").Value = "'98.23"
$ws.Range("E6").Value = "  +4.27%  "

$ws.Range("D7").Value = "'0.506"
$ws.Range("E7").Value = "  +1.73%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("E9").Value = "  +3.34%  "

$ws.Range("D10").Value = "'34.11"
$ws.Range("E10").Value = "  +3.96%  "

$ws.Range("E11").Value = "  +1.78%  "

$ws.Range("D12").Value = "'49.10"
$ws.Range("E12").Value = "  +1.75%  "

$ws.Range("E13").Value = "  +4.17%  "

$ws.Range("D14").Value = "'17.87"
$ws.Range("E14").Value = "  +16.43%  "

$ws.Range("D15").Value = "'6.80"
$ws.Range("E15").Value = "  +2.70%  "

$ws.Range("D16").Value = "2.662.02"

$ws.Range("D17").Value = "2.277.66"
$ws.Range("E17").Value = "  +0.81%  "

$ws.Range("D18").Value = "'0.808"
$ws.Range("E18").Value = "  +4.47%  "

$ws.Range("D19").Value = "42.884.25"
$ws.Range("E19").Value = "  +1.92%  "

$ws.Range("D20").Value = "'12.34"
$ws.Range("E20").Value = "  +8.61%  "

$ws.Range("D21").Value = "0.0₃0906"
$ws.Range("E21").Value = "  +2.18%  "

$ws.Range("D22").Value = "'6.11"
$ws.Range("E22").Value = "  +1.97%  "

$ws.Range("D23").Value = "'67.82"
$ws.Range("E23").Value = "  +2.23%  "

$ws.Range("D24").Value = "'237.43"
$ws.Range("E24").Value = "  +2.08%  "

$ws.Range("D25").Value = "'2.09"
$ws.Range("E25").Value = "  +9.19%  "

$ws.Range("E27").Value = "  +0.33%  "

$ws.Range("D28").Value = "'24.63"
$ws.Range("E28").Value = "  +3.74%  "

$ws.Range("D29").Value = "'2.28"
$ws.Range("E29").Value = "  +11.76%  "

$ws.Range("D30").Value = "'166.10"
$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("D31").Value = "'33.92"
$ws.Range("E31").Value = "  +1.11%  "

$ws.Range("D32").Value = "'9.18"
$ws.Range("E32").Value = "  +1.85%  "

$ws.Range("E33").Value = "  -0.06%  "

$ws.Range("D34").Value = "'5.02"
$ws.Range("E34").Value = "  +1.96%  "

$ws.Range("E35").Value = "  +4.40%  "

$ws.Range("D36").Value = "'4.56"
$ws.Range("E36").Value = "  +5.52%  "

$ws.Range("D37").Value = "'17.11"
$ws.Range("E37").Value = "  +7.33%  "

$ws.Range("D38").Value = "'0.0700"
$ws.Range("E38").Value = "  +1.36%  "

$ws.Range("E39").Value = "  +3.85%  "

$ws.Range("E40").Value = "  +1.32%  "

$ws.Range("E41").Value = "  +5.01%  "

$ws.Range("D42").Value = "'0.110"
$ws.Range("E42").Value = "  +0.42%  "

$ws.Range("D44").Value = "1.996.96"
$ws.Range("E44").Value = "  +3.25%  "

$ws.Range("E45").Value = "  +2.92%  "

$ws.Range("D46").Value = "'10.05"
$ws.Range("E46").Value = "  +5.71%  "

$ws.Range("D47").Value = "'17.75"
$ws.Range("E47").Value = "  +2.41%  "

$ws.Range("D48").Value = "'2.87"
$ws.Range("E48").Value = "  +3.97%  "

$ws.Range("D49").Value = "'54.46"
$ws.Range("E49").Value = "  +4.75%  "

$ws.Range("D50").Value = "2.527.74"
$ws.Range("E50").Value = "  +1.84%  "

$ws.Range("E51").Value = "  +2.99%  "
